$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3049
$ws1.Range("F3").Value = 472
$ws1.Range("F4").Value = 51
$ws1.Range("F5").Value = 36
$ws1.Range("F7").Value = 1030
$ws1.Range("F8").Value = 14653
$ws1.Range("F9").Value = 171
$ws1.Range("F10").Value = 131
$ws1.Range("F11").Value = 5851
$ws1.Range("F12").Value = 596
$ws1.Range("F13").Value = 80
$ws1.Range("F14").Value = 46
$ws1.Range("F15").Value = 71
$ws1.Range("F17").Value = 17
$ws1.Range("F18").Value = 88
$ws1.Range("F19").Value = 187
$ws1.Range("F22").Value = 78
$ws1.Range("F23").Value = 10636
$ws1.Range("F24").Value = 1202
$ws1.Range("F25").Value = 67
$ws1.Range("F26").Value = 95
$ws1.Range("F27").Value = 3742

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3049
$ws4.Range("F4").Value = 472
$ws4.Range("F5").Value = 51
$ws4.Range("F6").Value = 36
$ws4.Range("F8").Value = 1030
$ws4.Range("F9").Value = 14653
$ws4.Range("F10").Value = 171
$ws4.Range("F11").Value = 131
$ws4.Range("F12").Value = 5851
$ws4.Range("F13").Value = 596
$ws4.Range("F14").Value = 80
$ws4.Range("F15").Value = 46
$ws4.Range("F16").Value = 71
$ws4.Range("F18").Value = 17
$ws4.Range("F19").Value = 88
$ws4.Range("F20").Value = 187
$ws4.Range("F23").Value = 78
$ws4.Range("F25").Value = 10636
$ws4.Range("F26").Value = 1202
$ws4.Range("F27").Value = 67
$ws4.Range("F28").Value = 95
$ws4.Range("F29").Value = 3742
